$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.100337
$ws.Range("H2").Value = 0.301011
$ws.Range("M2").Value = 33.16156633333333
$ws.Range("N2").Value = 99.48469900000001
$ws.Range("O2").Value = 0.2666406279966088
$ws.Range("P2").Value = 0.2666406279966088
$ws.Range("Q2").Value = 3.327332081187667
$ws.Range("R2").Value = 29.945988730689
$ws.Range("S2").Value = 0.2666406279966088
$ws.Range("T2").Value = 0.2666406279966088

# Row 3
$ws.Range("G3").Value = 0.100337
$ws.Range("H3").Value = 0.301011
$ws.Range("O3").Value = 0.2064675829184841
$ws.Range("P3").Value = 0.2064675829184841
$ws.Range("Q3").Value = 2.576449873868
$ws.Range("R3").Value = 23.188048864812
$ws.Range("S3").Value = 0.2064675829184841
$ws.Range("T3").Value = 0.2064675829184841

# Row 4
$ws.Range("G4").Value = 0.100337
$ws.Range("H4").Value = 0.301011
$ws.Range("M4").Value = 15.57007066666667
$ws.Range("N4").Value = 46.710212
$ws.Range("O4").Value = 0.1251935261073135
$ws.Range("P4").Value = 0.1251935261073135
$ws.Range("Q4").Value = 1.562254180481333
$ws.Range("R4").Value = 14.060287624332
$ws.Range("S4").Value = 0.1251935261073135
$ws.Range("T4").Value = 0.1251935261073135

# Row 5
$ws.Range("G5").Value = 0.100337
$ws.Range("H5").Value = 0.301011
$ws.Range("M5").Value = 49.95841666666666
$ws.Range("N5").Value = 149.87525
$ws.Range("O5").Value = 0.4016982629775935
$ws.Range("P5").Value = 0.4016982629775935
$ws.Range("Q5").Value = 5.012677653083333
$ws.Range("R5").Value = 45.11409887775
$ws.Range("S5").Value = 0.4016982629775935
$ws.Range("T5").Value = 0.4016982629775935
